$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing the existing rows 104-174 down to 105-175.
$ws.Rows("104:104").Insert()

# Populate the freshly inserted row 104 with the new record.
$ws.Range("A104").Value = 5
$ws.Range("B104").Value = "Macroferia Regional de Talca"
$ws.Range("C104").Value = "Maule"
$ws.Range("D104").Value = 45086
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100108
$ws.Range("H104").Value = "Tropicales y subtropicales"
$ws.Range("I104").Value = 100108002
$ws.Range("J104").Value = "Mango"
$ws.Range("K104").Value = "Sin especificar"
$ws.Range("L104").Value = "Primera"
$ws.Range("M104").Value = 250
$ws.Range("N104").Value = 7000
$ws.Range("O104").Value = 8000
$ws.Range("P104").Value = 7200
$ws.Range("Q104").Value = "$/bandeja 4 kilos"
$ws.Range("R104").Value = "Ecuador"
$ws.Range("S104").Value = 1800
$ws.Range("T104").Value = 4
